# Update Name of Algo
# Apply updated imputed values in column C for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = -12.20379999999999
    10  = -13.2484
    12  = -10.73879999999999
    18  = -12.54769999999999
    25  = -13.3178
    37  = -13.7472
    55  = -13.46669999999999
    68  = -11.5134
    77  = -12.50300000000001
    78  = -12.68430000000001
    79  = -12.1874
    80  = -13.23480000000001
    81  = -12.8191
    82  = -12.258
    84  = -13.18079999999999
    101 = -12.53750000000001
    102 = -13.2619
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
